$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 8.25
$ws.Range("O2").Value = 49.46
$ws.Range("Q2").Value = -7.25
$ws.Range("AJ2").Value = 8.33
$ws.Range("AK2").Value = 8.220000000000001
$ws.Range("AN2").Value = 8.33
$ws.Range("AO2").Value = 8.220000000000001

$ws.Range("H3").Value = 10.72
$ws.Range("O3").Value = 56.5
$ws.Range("P3").Value = 11.55
$ws.Range("Q3").Value = 7.79
$ws.Range("AJ3").Value = 10.9
$ws.Range("AK3").Value = 10.69
$ws.Range("AN3").Value = 10.9
$ws.Range("AO3").Value = 10.69

$ws.Range("H4").Value = 3.75
$ws.Range("O4").Value = 292.67
$ws.Range("Q4").Value = -56.31
$ws.Range("T4").Value = 3.93
$ws.Range("X4").Value = 3.93
$ws.Range("AB4").Value = 3.93
$ws.Range("AF4").Value = 3.93
$ws.Range("AJ4").Value = 3.75
$ws.Range("AN4").Value = 3.75

$ws.Range("H5").Value = 2
$ws.Range("O5").Value = 158.06
$ws.Range("Q5").Value = -25.56
$ws.Range("T5").Value = 2.03
$ws.Range("X5").Value = 2.03
$ws.Range("AB5").Value = 2.03
$ws.Range("AF5").Value = 2.03
$ws.Range("AJ5").Value = 2
$ws.Range("AK5").Value = 1.92
$ws.Range("AN5").Value = 2
$ws.Range("AO5").Value = 1.92

$ws.Range("H6").Value = 3.5
$ws.Range("O6").Value = 45.23
$ws.Range("Q6").Value = -7.35
